# "#5: cash & deposit done"
#
# The 存款 (deposits) sheet is expanded: the B:F columns (which previously
# mixed up bank/deposit-type/currency data into a mislabeled header row)
# get a correct header, and new columns G:M are added carrying the same
# property_category / category / date / legislator_name / legislator_id /
# source_file / index metadata that the other sheets in this workbook
# already have.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---------------------------------------------------------------------
# Header row (row 1)
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 2).Value = "bank"               # B1
$ws.Cells.Item(1, 3).Value = "deposit_type"       # C1
$ws.Cells.Item(1, 4).Value = "currency"           # D1
$ws.Cells.Item(1, 5).Value = "owner"              # E1
$ws.Cells.Item(1, 6).Value = "total"              # F1
$ws.Cells.Item(1, 7).Value = "property_category"  # G1
$ws.Cells.Item(1, 8).Value = "category"           # H1
$ws.Cells.Item(1, 9).Value = "date"               # I1
$ws.Cells.Item(1, 10).Value = "legislator_name"   # J1
$ws.Cells.Item(1, 11).Value = "legislator_id"     # K1
$ws.Cells.Item(1, 12).Value = "source_file"       # L1
$ws.Cells.Item(1, 13).Value = "index"             # M1

# Copy the existing (bold + bordered) header style onto the new header
# cells G1:M1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# ---------------------------------------------------------------------
# Data rows (rows 2-10), columns A-M
# ---------------------------------------------------------------------
# A bank               deposit_type  currency  owner    total    property_category category date        legislator_name legislator_id source_file index
$rows = @(
    @(49, "彰化商業銀行大直分行",             "活期儲蓄存款", "新臺幣", "王金平",   2202742, "deposit", "normal", "2012-03-20", "王金平", 22, "tmpd1a31", 49),
    @(50, "彰化商業銀行大直分行",             "定期存款",     "新臺幣", "王金平",   1000000, "deposit", "normal", "2012-03-20", "王金平", 22, "tmpd1a31", 50),
    @(51, "臺灣銀行群賢分行",                 "活期儲蓄存款", "新臺幣", "王金平",   2243195, "deposit", "normal", "2012-03-20", "王金平", 22, "tmpd1a31", 51),
    @(52, "臺灣銀行",                         "活期存款",     "新臺幣", "王陳彩蓮", 32888,   "deposit", "normal", "2012-03-20", "王金平", 22, "tmpd1a31", 52),
    @(53, "彰化商業銀行大直分行",             "活期儲蓄存款", "新臺幣", "王陳彩蓮", 1738052, "deposit", "normal", "2012-03-20", "王金平", 22, "tmpd1a31", 53),
    @(54, "彰化商業銀行大直分行",             "定期存款",     "新臺幣", "王陳彩蓮", 500000,  "deposit", "normal", "2012-03-20", "王金平", 22, "tmpd1a31", 54),
    @(55, "國泰世華商業銀行館前分行",         "活期儲蓄存款", "新臺幣", "王陳彩蓮", 1045697, "deposit", "normal", "2012-03-20", "王金平", 22, "tmpd1a31", 55),
    @(56, "台北富邦商業銀行和平分行",         "活期儲蓄存款", "新臺幣", "王陳彩蓮", 14025,   "deposit", "normal", "2012-03-20", "王金平", 22, "tmpd1a31", 56),
    @(57, "中華郵政股份有限公司台北104支郵局", "活期存款",    "新臺幣", "王陳彩蓮", 609811,  "deposit", "normal", "2012-03-20", "王金平", 22, "tmpd1a31", 57)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]    # A  index (kept as before)
    $ws.Cells.Item($r, 2).Value = $row[1]    # B  bank
    $ws.Cells.Item($r, 3).Value = $row[2]    # C  deposit_type
    $ws.Cells.Item($r, 4).Value = $row[3]    # D  currency
    $ws.Cells.Item($r, 5).Value = $row[4]    # E  owner
    $ws.Cells.Item($r, 6).Value = $row[5]    # F  total
    $ws.Cells.Item($r, 7).Value = $row[6]    # G  property_category
    $ws.Cells.Item($r, 8).Value = $row[7]    # H  category
    # I (date) handled separately below so Excel doesn't reinterpret the
    # "2012-03-20" text as a date serial number.
    $ws.Cells.Item($r, 10).Value = $row[9]   # J  legislator_name
    $ws.Cells.Item($r, 11).Value = $row[10]  # K  legislator_id
    $ws.Cells.Item($r, 12).Value = $row[11]  # L  source_file
    $ws.Cells.Item($r, 13).Value = $row[12]  # M  index
}

# Copy the existing (plain) data-row style onto the new data cells
# G2:M10 in one shot.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("G2:M10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Now fill in the date column as literal text (not an Excel date value).
$ws.Range("I2:I10").NumberFormat = "@"
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 9).Value = $rows[$i][8]  # I date
}
